$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.188.48'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '3.143.89'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '242.48'
$ws.Range('E5').Value = '  +2.11%  '
$ws.Range('D6').Value = '627.55'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').Value = '1.16'
$ws.Range('E7').Value = '  +10.78%  '
$ws.Range('D8').Value = '0.374'
$ws.Range('E8').Value = '  +5.27%  '
$ws.Range('D10').Value = '3.140.86'
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('E11').Value = '  +6.11%  '
$ws.Range('D12').Value = '0.206'
$ws.Range('E12').Value = '  +4.30%  '
$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  +4.63%  '
$ws.Range('D14').Value = '35.93'
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('D15').Value = '5.54'
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('D16').Value = '90.711.86'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '3.725.27'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').Value = '3.141.33'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').Value = '3.80'
$ws.Range('E19').Value = '  +3.00%  '
$ws.Range('D20').Value = '14.76'
$ws.Range('E20').Value = '  +2.13%  '
$ws.Range('B21').Value = 'PEPE'
$ws.Range('C21').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D21').Value = '0.0000215'
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').Value = '5.91'
$ws.Range('E22').Value = '  +4.08%  '
$ws.Range('D23').Value = '452.65'
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('D24').Value = '9.18'
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('D25').Value = '6.00'
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('D26').Value = '93.72'
$ws.Range('E26').Value = '  +3.94%  '
$ws.Range('D27').Value = '12.03'
$ws.Range('E27').Value = '  -3.96%  '
$ws.Range('D28').Value = '3.284.76'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '0.181'
$ws.Range('E30').Value = '  +13.18%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.123'
$ws.Range('E31').Value = '  +43.81%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.227'
$ws.Range('E32').Value = '  +13.28%  '
$ws.Range('D33').Value = '9.16'
$ws.Range('E33').Value = '  -4.02%  '
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +25.49%  '
$ws.Range('D35').Value = '0.164'
$ws.Range('E35').Value = '  +9.27%  '
$ws.Range('D36').Value = '27.01'
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('D37').Value = '7.74'
$ws.Range('E37').Value = '  +8.92%  '
$ws.Range('D38').Value = '4.18'
$ws.Range('E38').Value = '  +23.48%  '
$ws.Range('D39').Value = '503.22'
$ws.Range('E39').Value = '  -1.68%  '
$ws.Range('D40').Value = '1.94'
$ws.Range('E40').Value = '  +0.64%  '
$ws.Range('D41').Value = '3.66'
$ws.Range('E41').Value = '  -4.07%  '
$ws.Range('D42').Value = '1.31'
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('D43').Value = '0.428'
$ws.Range('E43').Value = '  +0.46%  '
$ws.Range('D44').Value = '22.13'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').Value = '1.94'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '156.55'
$ws.Range('E47').Value = '  +4.78%  '
$ws.Range('D48').Value = '0.704'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('D49').Value = '4.60'
$ws.Range('E49').Value = '  +1.45%  '
$ws.Range('D50').Value = '1.36'
$ws.Range('E50').Value = '  +0.82%  '
$ws.Range('D51').Value = '45.15'
$ws.Range('E51').Value = '  -1.05%  '
